# Apply updated cryptocurrency price/volume figures to the "cryptos" sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'28.478.61"
$ws.Range("E2").Value = "  -0.45%  "
$ws.Range("D3").Value = "'1.819.61"
$ws.Range("E3").Value = "  -0.43%  "
$ws.Range("D4").Value = "'1.002"
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").Value = "'316.88"
$ws.Range("E5").Value = "  +0.05%  "
$ws.Range("E6").Value = "  +0.07%  "
$ws.Range("D7").Value = "'0.5164"
$ws.Range("E7").Value = "  -3.08%  "
$ws.Range("D8").Value = "'0.3888"
$ws.Range("E8").Value = "  -2.10%  "
$ws.Range("D9").Value = "'0.08468"
$ws.Range("E9").Value = "  +9.05%  "
$ws.Range("D10").Value = "'41.81"
$ws.Range("E10").Value = "  -0.57%  "
$ws.Range("E11").Value = "  -0.73%  "
$ws.Range("D12").Value = "'6.434"
$ws.Range("E12").Value = "  +1.78%  "
$ws.Range("E13").Value = "  -0.54%  "
$ws.Range("D14").Value = "'1.002"
$ws.Range("E14").Value = "  +0.03%  "
$ws.Range("D15").Value = "'7.512"
$ws.Range("E15").Value = "  -0.80%  "
$ws.Range("D16").Value = "'1.818.43"
$ws.Range("E16").Value = "  -0.02%  "
$ws.Range("D17").Value = "'0.00001142"
$ws.Range("E17").Value = "  +4.71%  "
$ws.Range("D18").Value = "'92.80"
$ws.Range("E18").Value = "  -0.78%  "
$ws.Range("D19").Value = "'0.06645"
$ws.Range("E19").Value = "  +0.40%  "
$ws.Range("D20").Value = "'17.75"
$ws.Range("E20").Value = "  -0.22%  "
$ws.Range("E21").Value = "  +0.07%  "
$ws.Range("D22").Value = "'6.089"
$ws.Range("E22").Value = "  +0.11%  "
$ws.Range("D23").Value = "'28.511.57"
$ws.Range("E23").Value = "  -0.34%  "
$ws.Range("D24").Value = "'11.45"
$ws.Range("E24").Value = "  +2.21%  "
$ws.Range("D25").Value = "'2.272"
$ws.Range("E25").Value = "  +1.64%  "
$ws.Range("D26").Value = "'21.03"
$ws.Range("E26").Value = "  +1.08%  "
$ws.Range("E27").Value = "  +1.53%  "
$ws.Range("D28").Value = "'2.029.99"
$ws.Range("E28").Value = "  -0.04%  "
$ws.Range("D29").Value = "'2.398"
$ws.Range("E29").Value = "  -1.01%  "
$ws.Range("D30").Value = "'125.64"
$ws.Range("E30").Value = "  +0.31%  "
$ws.Range("D31").Value = "'0.1089"
$ws.Range("E31").Value = "  -3.55%  "
$ws.Range("D32").Value = "'1.095"
$ws.Range("E32").Value = "  -4.97%  "
$ws.Range("D33").Value = "'5.724"
$ws.Range("E33").Value = "  -0.44%  "
$ws.Range("D34").Value = "'0.07441"
$ws.Range("E34").Value = "  +1.42%  "
$ws.Range("D35").Value = "'3.652"
$ws.Range("E35").Value = "  -0.31%  "
$ws.Range("D36").Value = "'0.2234"
$ws.Range("E36").Value = "  -1.71%  "
$ws.Range("D37").Value = "'0.02360"
$ws.Range("E37").Value = "  +0.33%  "
$ws.Range("D38").Value = "'5.206"
$ws.Range("E38").Value = "  -0.08%  "
$ws.Range("D39").Value = "'8.844"
$ws.Range("E39").Value = "  -0.81%  "
$ws.Range("D40").Value = "'0.6319"
$ws.Range("E40").Value = "  +0.19%  "
$ws.Range("D41").Value = "'11.26"
$ws.Range("E41").Value = "  -1.47%  "
$ws.Range("D42").Value = "'1.194"
$ws.Range("E42").Value = "  -0.22%  "
$ws.Range("D43").Value = "'1.401"
$ws.Range("E43").Value = "  +0.15%  "
$ws.Range("D44").Value = "'13.56"
$ws.Range("E44").Value = "  +0.03%  "
$ws.Range("D45").Value = "'3.782"
$ws.Range("E45").Value = "  +1.61%  "
$ws.Range("D46").Value = "'0.5945"
$ws.Range("E46").Value = "  -0.09%  "
$ws.Range("D47").Value = "'126.25"
$ws.Range("E47").Value = "  +0.52%  "
$ws.Range("E48").Value = "  -0.64%  "
$ws.Range("D49").Value = "'1.203"
$ws.Range("E49").Value = "  +0.78%  "
$ws.Range("D50").Value = "'0.06973"
$ws.Range("E50").Value = "  +0.13%  "
$ws.Range("D51").Value = "'74.35"
$ws.Range("E51").Value = "  -0.23%  "
